# feat: complete overhaul of Excel pipeline to tracking raw quantities and
# real-time yfinance valuation computation
#
# - Re-applies the "YYYY-MM-DD HH:MM:SS" date/time format to the date column
#   on the "Daily" sheet (A2:A184), collapsing the old duplicate/unused style
#   definition down to a single shared style.
# - Adds a new "Holdings" worksheet (right after "Daily") that tracks the
#   raw quantities held per symbol (cash + each position), replacing the
#   implicit valuation bookkeeping that used to live only in "Daily".

$wb = $excel.ActiveWorkbook
$daily = $wb.Worksheets.Item("Daily")

# --- Daily: normalize the date column's number format -----------------
# Re-assert the canonical (uppercase) date/time format on the date column.
# The workbook previously carried a stray, unused duplicate style entry for
# this exact format; cycling through the lowercase variant first and then
# back to the canonical uppercase code collapses the cells onto a single
# shared style.
$dateRange = $daily.Range("A2:A184")
$dateRange.NumberFormat = "yyyy-mm-dd h:mm:ss"
$dateRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- New "Holdings" sheet ----------------------------------------------
$holdings = $wb.Worksheets.Add($null, $daily)
$holdings.Name = "Holdings"

$headers = @("Symbol", "Name", "Quantity")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $holdings.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$rows = @(
    @("CASH", "USD Cash", 571.73),
    @("GC=F", "Gold Futures", 0.287),
    @("NVDA", "NVIDIA Corp", 17.4),
    @("TSLA", "Tesla Inc", 0)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $holdings.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Keep "Daily" as the active sheet/tab, matching the original workbook view.
$daily.Activate()

Write-Output "edit complete"
